# "Generate Report for Handoff"
# The localization file fa92f36a-768f-4af6-8053-a628f46e61f8.md moved from
# "In Translation" to "Ready for handoff" for both the zh-cn and de-de
# locales, with refreshed handoff timestamps, and the zh-cn row's Priority
# flipped from "ht" (human translation) to "mt" (machine translation).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-28 20:15:19"
$ov.Columns.Item(5).ColumnWidth = 16.3
$ov.Columns.Item(6).ColumnWidth = 16.3

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-08-28 20:15:15"
$zh.Columns.Item(3).ColumnWidth = 16.3

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("H3").Value = "2016-08-28 20:15:19"
$de.Columns.Item(3).ColumnWidth = 16.3
